# Weekly update: insert a new latest-week record at the top of the data
# block (row 594) for "Hortaliza, Terminal Hortofrutícola Agro Chillán -
# Brócoli". All existing data rows 594:638 shift down by one to 595:639.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 594, pushing the existing rows 594-638 down
# to 595-639 (this also grows the sheet dimension to A1:R639 automatically).
$ws.Rows("594").Insert()

# Populate the new row 594 with the new weekly record.
$ws.Cells.Item(594, 1).Value  = 7
$ws.Cells.Item(594, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(594, 3).Value  = "Ñuble"
$ws.Cells.Item(594, 4).Value  = 45223
$ws.Cells.Item(594, 5).Value  = 16
$ws.Cells.Item(594, 6).Value  = 100112023
$ws.Cells.Item(594, 7).Value  = "Brócoli"
$ws.Cells.Item(594, 8).Value  = "Sin especificar"
$ws.Cells.Item(594, 9).Value  = "Primera"
$ws.Cells.Item(594, 10).Value = 400
$ws.Cells.Item(594, 11).Value = 1200
$ws.Cells.Item(594, 12).Value = 1200
$ws.Cells.Item(594, 13).Value = 1200
$ws.Cells.Item(594, 14).Value = "$/unidad"
$ws.Cells.Item(594, 15).Value = "Región del Maule"
$ws.Cells.Item(594, 16).Value = 1200
$ws.Cells.Item(594, 17).Value = 1
$ws.Cells.Item(594, 18).Value = "Hortaliza"
